$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Fgf1", "Fgfr2", "ECs", 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 1.684271952582, 15.158447573238, 0.03953876810189297, 0.03953876810189298)
    ,@("ECs", "Fgf1", "Fgfr2", "FAPs", 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 2.783107476884889, 25.047967291964, 0.06533424781105156, 0.06533424781105156)
    ,@("ECs", "Fgf1", "Fgfr2", "MuSCs", 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.3740528624806667, 3.366475762326, 0.008780998439592727, 0.008780998439592728)
    ,@("FAPs", "Fgf1", "Fgfr2", "ECs", 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 7.310862291651, 65.797760624859, 0.1716245932441777, 0.1716245932441777)
    ,@("FAPs", "Fgf1", "Fgfr2", "FAPs", 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 12.08054048230022, 108.724864340702, 0.2835941594485134, 0.2835941594485134)
    ,@("FAPs", "Fgf1", "Fgfr2", "MuSCs", 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 1.623638607293667, 14.612747465643, 0.03811538289683612, 0.03811538289683613)
    ,@("MuSCs", "Fgf1", "Fgfr2", "ECs", 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 5.710986372042001, 51.398877348378, 0.1340670462695045, 0.1340670462695045)
    ,@("MuSCs", "Fgf1", "Fgfr2", "FAPs", 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 9.436889837209334, 84.93200853488399, 0.2215337008400159, 0.2215337008400159)
    ,@("MuSCs", "Fgf1", "Fgfr2", "MuSCs", 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 1.268328904234, 11.414960138106, 0.02977438551093763, 0.02977438551093764)
    ,@("Resolving-Mac", "Fgf1", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 0.113170740717, 1.018536666453, 0.002656715660597019, 0.002656715660597019)
    ,@("Resolving-Mac", "Fgf1", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 0.1870044408037778, 1.683039967234, 0.004389982987978754, 0.004389982987978755)
    ,@("Resolving-Mac", "Fgf1", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.02513361304233333, 0.226202517381, 0.0005900187889016981, 0.0005900187889016983)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
